# Apply a green highlight to the paragraph "שינויי שם קובץ" (file name
# change) — both to the run text itself and to the paragraph mark, as in
# the source diff.
#
# wdColorIndex: wdBrightGreen = 4  ->  <w:highlight w:val="green"/>

$d = $word.ActiveDocument
$targetText = "שינויי שם קובץ"

$rng = $d.Content
$rng.Find.Execute($targetText, $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    # Cut the run out of the paragraph so it becomes run-less; applying
    # HighlightColorIndex to a run-less paragraph range stamps the
    # paragraph-mark run properties (w:pPr/w:rPr) instead of a run's
    # w:rPr. Then paste the exact same (fully formatted) run back.
    $rng.Cut()
    $rng.HighlightColorIndex = 4
    $rng.Paste()

    # Re-locate the just-restored run and highlight it too, so the
    # visible text itself is highlighted like in the source edit.
    $rng2 = $d.Content
    $rng2.Find.Execute($targetText, $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
    if ($rng2.Find.Found) {
        $rng2.HighlightColorIndex = 4
    }
}
